# Auto update stock data
# Update Date_1 (column A) and EBITDA (column B) for the latest data rows.
# Each row gets a refreshed date; most rows also get an updated EBITDA figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    # Force text entry so date-like / numeric-like strings aren't
    # auto-converted to a real date serial or number by the COM layer.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    # Restore the cell to its original (unstyled) formatting state.
    $rng.ClearFormats()
}

# Row 2 - Alcoa
Set-TextValue "A2" "2026/01/09"
Set-TextValue "B2" "7.63"

# Row 8 - Rio Tinto
Set-TextValue "A8" "2026/01/09"
Set-TextValue "B8" "8.78"

# Row 14 - Norsk Hydro
Set-TextValue "A14" "2026/01/09"
Set-TextValue "B14" "3.09"

# Row 20 - Reliance
Set-TextValue "A20" "2026/01/09"
Set-TextValue "B20" "13.45"

# Row 26 - Kaiser
Set-TextValue "A26" "2026/01/09"
Set-TextValue "B26" "11.61"

# Row 32 - Ryerson
Set-TextValue "A32" "2026/01/09"
Set-TextValue "B32" "28.51"

# Row 38 - Alro Steel (date only, EBITDA unchanged)
Set-TextValue "A38" "2026/01/09"

# Row 44 - Ultra
Set-TextValue "A44" "2026/01/09"
Set-TextValue "B44" "13.39"

# Row 50 - Benchmark
Set-TextValue "A50" "2026/01/09"
Set-TextValue "B50" "11.50"

# Row 56 - Celestica
Set-TextValue "A56" "2026/01/09"
Set-TextValue "B56" "30.23"

# Row 62 - Jabil
Set-TextValue "A62" "2026/01/09"
Set-TextValue "B62" "11.10"

# Row 68 - Flex
Set-TextValue "A68" "2026/01/09"
Set-TextValue "B68" "12.61"

# Row 74 - MKS
Set-TextValue "A74" "2026/01/09"
Set-TextValue "B74" "17.88"
